# [DOM-306] Finish acceptance csenarios
# Update the scenario list for "auth_13" (row 21) from a 4-item list
# to the finished 2-item list, and scroll/select to reflect where the
# author was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = "ТС-НАСТР1, ТС-НАСТР2"

$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D20").Select()
